# Milestone3Results.xlsx - "Add files via upload" edit
# Renames Sheet1 -> RandomForest, and fills in a RandomForest accuracy-vs-j
# results table (B2:H26) including a Lucida-Console styled block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet -------------------------------------------------
$ws.Name = "RandomForest"

# --- B2:B12 - blank cells carrying the "console output" font ----------
# Build the format once on a scratch cell, then copy/paste-special the
# format onto the target range so the style table only grows by the one
# style actually used (matches how Excel collapses style churn).
$helper = $ws.Range("Z1")
$helper.Font.Size = 10
$helper.Font.Color = 526344        # RGB(8,8,8) = FF080808
$helper.Font.Name = "Lucida Console"
$helper.VerticalAlignment = -4108  # -4108 = xlVAlignCenter

$helper.Copy()
$ws.Range("B2:B12").PasteSpecial(-4122)  # -4122 = xlPasteFormats
$helper.Clear()

# --- Header row 16 ------------------------------------------------------
$ws.Range("H16").Value = "j"
$ws.Range("D16").Value = "genTrainAcc"
$ws.Range("E16").Value = "0TrainAcc"
$ws.Range("F16").Value = "genTestAcc"
$ws.Range("G16").Value = "0TestAcc"

# --- Data rows 17-26 ------------------------------------------------------
$ws.Range("C17").Value = "accRF.rf..test..train.classRebalance.1..test..train....."

$data = @(
  @(17, 0.97188085205992503, 0.97313904494381998, 0.48167950021792799, 0.95884553714591103, 1),
  @(18, 0.96972534332084903, 0.949730805243446,   0.55150370477989297, 0.93907001603420603, 2),
  @(19, 0.97157361891385796, 0.92532771535580505, 0.60177248292895502, 0.91608765366114397, 3),
  @(20, 0.97257724719101102, 0.89981273408239704, 0.64564869969490002, 0.89791555318011795, 4),
  @(21, 0.97348041510611705, 0.88541666666666696, 0.69173325584774104, 0.87814003206841296, 5),
  @(22, 0.97447665864098498, 0.86698267790262196, 0.71988958303065498, 0.86424371993586302, 6),
  @(23, 0.97505559456928803, 0.85059691011236005, 0.73627778584919401, 0.84981293425975402, 7),
  @(24, 0.97558390553474805, 0.83257256554307102, 0.75815778003777401, 0.82843399251737004, 8),
  @(25, 0.97559105805243496, 0.79950842696629199, 0.78727299142815599, 0.79903794762159297, 9),
  @(26, 0.97632043752127995, 0.79037921348314599, 0.79715240447479296, 0.79048637092463903, 10)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Range("D$r").Value = $row[1]
  $ws.Range("E$r").Value = $row[2]
  $ws.Range("F$r").Value = $row[3]
  $ws.Range("G$r").Value = $row[4]
  $ws.Range("H$r").Value = $row[5]
}

# --- Page setup -----------------------------------------------------------
$ws.PageSetup.Orientation = 1  # 1 = xlPortrait

# --- View selection --------------------------------------------------------
$ws.Range("I33").Select() | Out-Null

Write-Output "edit applied"
